# Auto-generated edit script: updates cryptos price/volume table cells
# to match the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.090.40"
$ws.Range("E2").Value = "  -4.29%  "
$ws.Range("D3").Value = "2.959.50"
$ws.Range("E3").Value = "  -6.65%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'571.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.24%  "
$ws.Range("D6").Value = "'123.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.61%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "2.950.01"
$ws.Range("E8").Value = "  -6.87%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.87%  "
$ws.Range("D10").Value = "'0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.76%  "
$ws.Range("D11").Value = "'5.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("D12").Value = "'0.433"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.40%  "
$ws.Range("D13").Value = "'0.0000221"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.37%  "
$ws.Range("D14").Value = "'32.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.44%  "
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "3.467.46"
$ws.Range("E16").Value = "  -6.14%  "
$ws.Range("D17").Value = "60.169.86"
$ws.Range("E17").Value = "  -4.15%  "
$ws.Range("D18").Value = "2.970.93"
$ws.Range("E18").Value = "  -6.33%  "
$ws.Range("D19").Value = "'6.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.41%  "
$ws.Range("D20").Value = "'426.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.43%  "
$ws.Range("D21").Value = "'12.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.26%  "
$ws.Range("D22").Value = "'0.656"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.39%  "
$ws.Range("D23").Value = "'6.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.75%  "
$ws.Range("D24").Value = "'12.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.72%  "
$ws.Range("D25").Value = "'78.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.70%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'2.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.68%  "
$ws.Range("D29").Value = "'7.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.27%  "
$ws.Range("D30").Value = "'1.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.24%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'6.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.58%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'25.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.67%  "
$ws.Range("D33").Value = "'0.0917"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.33%  "
$ws.Range("D34").Value = "'2.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.20%  "
$ws.Range("D35").Value = "'0.938"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.94%  "
$ws.Range("D36").Value = "'5.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.99%  "
$ws.Range("D37").Value = "'49.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.88%  "
$ws.Range("D38").Value = "0.0₃0645"
$ws.Range("E38").Value = "  -8.21%  "
$ws.Range("D39").Value = "'0.0355"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.47%  "
$ws.Range("D40").Value = "'7.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("D41").Value = "'0.108"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("D42").Value = "'375.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.91%  "
$ws.Range("D43").Value = "2.629.52"
$ws.Range("E43").Value = "  -5.61%  "
$ws.Range("D44").Value = "'2.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.16%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'0.233"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.06%  "
$ws.Range("D47").Value = "'119.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.32%  "
$ws.Range("D48").Value = "'1.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.35%  "
$ws.Range("D49").Value = "'0.105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.06%  "
$ws.Range("D50").Value = "'23.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.51%  "
$ws.Range("D51").Value = "'30.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -11.13%  "
